$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New gRNA analysis results changed a few existing counts.
$ws.Cells.Item(3, 2).Value = 57
$ws.Cells.Item(3, 3).Value = 32

$ws.Cells.Item(4, 2).Value = 148
$ws.Cells.Item(4, 3).Value = 89

# Insert a new "Tat " region row above the existing "Tat/Rev" row (row 7),
# pushing Tat/Rev, Vif and Vpr down by one row.
$ws.Rows.Item(7).Insert()

$cell = $ws.Cells.Item(7, 1)
$cell.Value = "Tat "
$cell.Font.Bold = $true
$cell.HorizontalAlignment = -4108
$cell.VerticalAlignment = -4160
$cell.Borders.LineStyle = 1

$ws.Cells.Item(7, 2).Value = 1
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = 0
